$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-report values between row 3 and row 4
# (dates D3/D4, volume J, prices K/L/M, and P column) to reflect
# the corrected weekly ordering.

# Row 3 new values (previously row 4's values)
$ws.Range("D3").Value = 44379
$ws.Range("J3").Value = 240
$ws.Range("K3").Value = 31000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 31500
$ws.Range("P3").Value = 1260

# Row 4 new values (previously row 3's values)
$ws.Range("D4").Value = 44827
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 31000
$ws.Range("M4").Value = 30500
$ws.Range("P4").Value = 1220
